$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (e.g. trailing zeros,
# thousand-dot formatting) instead of being auto-coerced to numbers/dates by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.592.96'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.000.55'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -4.16%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.012'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.82%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '330.01'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.86%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5004'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4223'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -4.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.84'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09011'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.118'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.36'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -5.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.015.60'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.067'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -6.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.468'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -6.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.013'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '94.08'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -6.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001112'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06658'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.958'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -6.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.605.82'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.97'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.304'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.78'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.68'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.426'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.298'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -8.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.28'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.052'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -6.82%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.578'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -5.25%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09935'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.03%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -6.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.799'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02469'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.318'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -7.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.310'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06350'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -6.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6566'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -5.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.68'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -6.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2052'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -7.21%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6341'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -6.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.43'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -6.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.200'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -5.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.305'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.509'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.39%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06983'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.125'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -7.19%  '
